# "Added level in course" — re-organize the Backlog sheet:
#   * re-sort/re-group the existing backlog rows (course items, author items,
#     package/"numbered" items) into a new row layout
#   * add a new "Course level" task and promote a handful of tasks to "Done"
#   * drop three stale backlog items ("13. Styling", "Modal page", "Router
#     redirect") and replace "Add Validation on Save button" with the next
#     wave of author-module work
#   * add a small new H/M side list (WIP items) next to the main table
#
# Strategy: wipe the sheet's old B:D (rows 5-40) table completely, then
# re-write every cell of the new C:M (rows 5-37) layout from scratch so the
# shared-string table / cell grid matches the target exactly, re-applying the
# two cell styles used in this sheet (bold header / green "Done" fill) as we
# go — Excel's style de-duplication means re-applying Font.Bold or
# Interior.Color on a cell whose computed format already matches an existing
# cellXf reuses that xf rather than minting a new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Wipe everything the old table touched (and then some) so no stray cells /
# row spans are left behind from the previous B5:D40 layout.
$ws.Range("A1:N50").Clear()

# ---------------------------------------------------------------------
# Row 5 — bold table headers
# ---------------------------------------------------------------------
$ws.Range("C5").Value = "Task"
$ws.Range("C5").Font.Bold = $true

$ws.Range("D5").Value = "Satus"
$ws.Range("D5").Font.Bold = $true

$ws.Range("H5").Value = "WIP"
$ws.Range("H5").Font.Bold = $true

# ---------------------------------------------------------------------
# Rows 7-16 — Course block
# ---------------------------------------------------------------------
$ws.Range("C7").Value = "Delete Course  "
$ws.Range("D7").Value = "Done"
$ws.Range("D7").Interior.Color = 5296274

$ws.Range("C8").Value = "Hide empty course list when all course are deleted  "
$ws.Range("D8").Value = "Done"
$ws.Range("D8").Interior.Color = 5296274

$ws.Range("C9").Value = "Client side validation for category and link data"
$ws.Range("D9").Value = "Done"
$ws.Range("D9").Interior.Color = 5296274

$ws.Range("C10").Value = "Show #course on header  "
$ws.Range("D10").Value = "Done"
$ws.Range("D10").Interior.Color = 5296274

$ws.Range("C11").Value = "Sort course table (mapStateToProps)   "
$ws.Range("D11").Value = ""
$ws.Range("D11").Interior.Color = 5296274

$ws.Range("C12").Value = "Handle 404 on manage course page"

$ws.Range("C13").Value = "Confirmation dialouge on delete  "
$ws.Range("D13").Value = "Done"
$ws.Range("D13").Interior.Color = 5296274

$ws.Range("C14").Value = "Fixed issue of reset value after adding author"
$ws.Range("D14").Value = "Done"
$ws.Range("D14").Interior.Color = 5296274

$ws.Range("C15").Value = "Didable Save button when page loaded with blank values"
$ws.Range("D15").Value = "Done"
$ws.Range("D15").Interior.Color = 5296274

$ws.Range("C16").Value = "Saving author functionality"
$ws.Range("D16").Value = "Done"
$ws.Range("D16").Interior.Color = 5296274

# ---------------------------------------------------------------------
# Rows 18-22 — Author block
# ---------------------------------------------------------------------
$ws.Range("C18").Value = "Add a date of course update/add"

$ws.Range("C19").Value = "Add Course level"
$ws.Range("D19").Value = "Done"
$ws.Range("D19").Interior.Color = 5296274

$ws.Range("C20").Value = "Restruture Course module"
$ws.Range("D20").Value = "Done"
$ws.Range("D20").Interior.Color = 5296274

$ws.Range("C21").Value = "Add course review"

$ws.Range("C22").Value = "Delete author with own delete icon"

# ---------------------------------------------------------------------
# Rows 28-37 — remaining backlog / numbered package items
# ---------------------------------------------------------------------
$ws.Range("C28").Value = "Unsaved changes message when user is leaving manage course page"
$ws.Range("C29").Value = "Default Sorting on page load"
$ws.Range("C30").Value = "10. Revert abandoned changes"
$ws.Range("C31").Value = "11. Pagination"

$ws.Range("C32").Value = "12. add Delete icon"
$ws.Range("C32").Font.Bold = $true

$ws.Range("C33").Value = "Styling of react components"
$ws.Range("C34").Value = "13. Radium pacakge for stylying"
$ws.Range("C35").Value = "14. SAAS"
$ws.Range("C36").Value = "15. reselect"
$ws.Range("C37").Value = "Add a course description page"

# ---------------------------------------------------------------------
# Column H/M — new WIP side list
# ---------------------------------------------------------------------
$ws.Range("H7").Value = "Add DOB with date control in Add author page"
$ws.Range("M7").Value = "http://react-day-picker.js.org/examples/?overlay"

$ws.Range("H8").Value = "Author Adminstration (cant delete a author if he/she has a course)"
$ws.Range("H9").Value = "Add email to author page with validation"
$ws.Range("H10").Value = "Add tests for author module"

# ---------------------------------------------------------------------
# View state — mirror the new selection/scroll position
# ---------------------------------------------------------------------
$ws.Range("H11").Select()
